$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/29/2025  Through  10/5/2025"

# --- Numeric cell updates ---
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -60
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -26.666666666666
$ws.Range("L15").Value = -21.428571428571
$ws.Range("N15").Value = -42.105263157894
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 110
$ws.Range("J16").Value = 107
$ws.Range("K16").Value = 2.803738317757
$ws.Range("L16").Value = -9.090909090909
$ws.Range("M16").Value = -6.779661016949
$ws.Range("N16").Value = -83.751846381093
$ws.Range("C17").Value = 9
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 220
$ws.Range("J17").Value = 214
$ws.Range("K17").Value = 2.803738317757
$ws.Range("L17").Value = 24.293785310734
$ws.Range("M17").Value = 175
$ws.Range("N17").Value = -10.931174089068
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = -3.333333333333
$ws.Range("L18").Value = -30.952380952381
$ws.Range("M18").Value = -43.414634146341
$ws.Range("N18").Value = -92.592592592592
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 32.5
$ws.Range("I19").Value = 454
$ws.Range("J19").Value = 467
$ws.Range("K19").Value = -2.783725910064
$ws.Range("L19").Value = -10.629921259842
$ws.Range("M19").Value = 57.638888888888
$ws.Range("N19").Value = -58.988256549232
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -42.857142857142
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -41.379310344827
$ws.Range("I20").Value = 175
$ws.Range("J20").Value = 195
$ws.Range("K20").Value = -10.256410256410
$ws.Range("L20").Value = 8.024691358024
$ws.Range("M20").Value = 21.527777777777
$ws.Range("N20").Value = -93.201243201243
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 47.058823529411
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -6.779661016949
$ws.Range("I21").Value = 1088
$ws.Range("J21").Value = 1119
$ws.Range("K21").Value = -2.770330652368
$ws.Range("L21").Value = -5.719237435008
$ws.Range("M21").Value = 29.216152019002
$ws.Range("N21").Value = -82.437449556093
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -25
$ws.Range("M22").Value = -14.285714285714
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 48
$ws.Range("K23").Value = -14.285714285714
$ws.Range("L23").Value = -2.040816326530
$ws.Range("M23").Value = 140
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 71.428571428571
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 25.925925925925
$ws.Range("I24").Value = 857
$ws.Range("J24").Value = 895
$ws.Range("K24").Value = -4.245810055865
$ws.Range("L24").Value = -12.012320328542
$ws.Range("M24").Value = 5.153374233128
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = -26.666666666666
$ws.Range("I25").Value = 319
$ws.Range("J25").Value = 385
$ws.Range("K25").Value = -17.142857142857
$ws.Range("L25").Value = -21.234567901234
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 10.256410256410
$ws.Range("I26").Value = 401
$ws.Range("J26").Value = 354
$ws.Range("K26").Value = 13.276836158192
$ws.Range("L26").Value = 32.781456953642
$ws.Range("M26").Value = 3.886010362694
$ws.Range("D27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -40
$ws.Range("L27").Value = -33.333333333333
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 6
$ws.Range("I28").Value = 41
$ws.Range("K28").Value = 64
$ws.Range("L28").Value = 24.242424242424
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 33.333333333333
$ws.Range("N29").Value = -87.096774193548
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = -83.333333333333
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 24

# --- Helper: write a text-valued cell while preserving the donor text style (s=13) ---
function Set-TextCell($ref, $text) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range("C14").Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

# --- Cells changing from numeric to text placeholders ("0" / "***.*") ---
Set-TextCell "F15" "0"
Set-TextCell "C16" "0"
Set-TextCell "D17" "0"
Set-TextCell "E17" "***.*"
Set-TextCell "F22" "0"
Set-TextCell "D23" "0"
Set-TextCell "E23" "***.*"
Set-TextCell "F27" "0"
Set-TextCell "G28" "0"
Set-TextCell "H28" "***.*"
Set-TextCell "D29" "0"
Set-TextCell "E29" "***.*"
Set-TextCell "D30" "0"
Set-TextCell "E30" "***.*"

$excel.CutCopyMode = 0
